# Fix column F (ASSISTS) so values are stored as real numbers instead of
# text, and correct the CHAMPION column (H) entries that were mis-tagged
# for rows where the player was actually playing Ahri the whole game.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F: convert the text "0"/"4"/"5"/"6" values (rows 2-41) into numeric values.
$assists = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 0;  6  = 0;  7  = 0;  8  = 0;  9  = 0; 10 = 0;
    11 = 0; 12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 0; 18 = 0; 19 = 0; 20 = 0;
    21 = 4; 22 = 4; 23 = 4; 24 = 4; 25 = 4; 26 = 4; 27 = 4;
    28 = 5; 29 = 5;
    30 = 6; 31 = 6; 32 = 6; 33 = 6; 34 = 6; 35 = 6; 36 = 6; 37 = 6; 38 = 6; 39 = 6; 40 = 6; 41 = 6
}

foreach ($row in $assists.Keys) {
    $ws.Cells.Item($row, 6).Value = $assists[$row]
}

# Column H: the champion these rows should show is "Ahri" (was wrongly
# recorded as Pyke / Gangplank / Taliyah / Akshan / Ornn).
$championRows = @(5, 11, 17, 23, 29, 30, 35, 37, 41)
foreach ($row in $championRows) {
    $ws.Cells.Item($row, 8).Value = "Ahri"
}
